$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the edited Price/Volume cells stay plain text (not auto-converted
# to numbers/percentages) by forcing a Text number format before writing,
# matching the original inlineStr string content in the sheet.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D2").Value = "323.93"
$ws.Range("E2").Value = "-2.58%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D3").Value = "44.53"
$ws.Range("E3").Value = "0.95%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.500"
$ws.Range("E4").Value = "-4.48%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08033"
$ws.Range("E5").Value = "-3.67%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D6").Value = "8.659"
$ws.Range("E6").Value = "-1.76%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D7").Value = "4.334"
$ws.Range("E7").Value = "-3.82%"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.13%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D9").Value = "2.690"
$ws.Range("E9").Value = "-7.03%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9408"
$ws.Range("E10").Value = "0.94%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1174"
$ws.Range("E11").Value = "-6.02%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1878"
$ws.Range("E12").Value = "-3.94%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09914"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04175"
$ws.Range("E14").Value = "5.51%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1065"
$ws.Range("E15").Value = "-0.15%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001271"
$ws.Range("E16").Value = "-2.54%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005920"
$ws.Range("E17").Value = "-0.24%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D18").Value = "3.594"
$ws.Range("E18").Value = "2.64%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D20").Value = "8.498"
$ws.Range("E20").Value = "-6.44%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1376"
$ws.Range("E21").Value = "1.02%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2646"
$ws.Range("E22").Value = "2.87%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04256"
$ws.Range("E23").Value = "-3.63%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001242"
$ws.Range("E24").Value = "-1.51%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004463"
$ws.Range("E25").Value = "2.21%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001236"
$ws.Range("E26").Value = "3.74%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.37%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02634"
$ws.Range("E39").Value = "-7.06%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05494"
$ws.Range("E40").Value = "-6.29%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007662"
$ws.Range("E41").Value = "-2.35%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.49%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.006960"
$ws.Range("E43").Value = "-23.37%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002056"
$ws.Range("E44").Value = "-2.20%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009214"
$ws.Range("E45").Value = "-12.05%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007129"
$ws.Range("E46").Value = "-1.87%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003416"
$ws.Range("E48").Value = "5.61%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.03%"
